# Update attendance_percentage, total_classes, and classes_attended
# for rows 2-6 on the active sheet (student_details data refresh /
# dynamic rendering update).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; C = 92.59999999999999; D = 47; E = 43 },
    @{ Row = 3; C = 87.8;              D = 40; E = 35 },
    @{ Row = 4; C = 60.8;              D = 40; E = 24 },
    @{ Row = 5; C = 93.7;              D = 45; E = 42 },
    @{ Row = 6; C = 72.7;              D = 43; E = 31 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 4).Value = $u.D
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
